$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update from the diff. Numeric-looking text values in column D
# (prices) are written with a leading apostrophe so Excel keeps them as literal
# text instead of silently parsing them into floating point numbers (which would
# lose formatting like trailing zeros). The style is reset to "Normal" afterward so
# the quote-prefix flag does not leave a stray style on the cell.

$ws.Range("D2").Value = "67.235.29"
$ws.Range("E2").Value = "  -3.36%  "
$ws.Range("D3").Value = "3.686.44"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'592.46"
$ws.Range("E5").Value = "  -3.51%  "
$ws.Range("D6").Value = "'167.03"
$ws.Range("E6").Value = "  -5.69%  "
$ws.Range("D7").Value = "3.685.93"
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.521"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -4.21%  "
$ws.Range("E11").Value = "  -4.34%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = "  -5.08%  "
$ws.Range("D13").Value = "'37.68"
$ws.Range("E13").Value = "  -5.39%  "
$ws.Range("D14").Value = "'0.0000240"
$ws.Range("E14").Value = "  -5.51%  "
$ws.Range("D15").Value = "4.295.58"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").Value = "3.679.27"
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").Value = "67.230.32"
$ws.Range("E17").Value = "  -3.40%  "
$ws.Range("E18").Value = "  -4.15%  "
$ws.Range("D19").Value = "'7.08"
$ws.Range("E19").Value = "  -6.16%  "
$ws.Range("D20").Value = "'16.83"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "'484.25"
$ws.Range("E21").Value = "  -4.85%  "
$ws.Range("D22").Value = "'9.08"
$ws.Range("E22").Value = "  -5.40%  "
$ws.Range("D23").Value = "'0.717"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("D24").Value = "'84.59"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "'2.30"
$ws.Range("E25").Value = "  -6.79%  "
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").Value = "'12.15"
$ws.Range("E27").Value = "  -5.40%  "
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").Value = "'9.96"
$ws.Range("E29").Value = "  -5.64%  "
$ws.Range("D30").Value = "'2.91"
$ws.Range("E30").Value = "  -2.61%  "
$ws.Range("D31").Value = "'2.34"
$ws.Range("E31").Value = "  -6.59%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'7.71"
$ws.Range("E32").Value = "  -5.14%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'31.70"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "3.819.28"
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("D35").Value = "3.620.89"
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = "  -6.86%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "'0.989"
$ws.Range("E38").Value = "  -5.36%  "
$ws.Range("D39").Value = "'5.75"
$ws.Range("E39").Value = "  -6.03%  "
$ws.Range("D40").Value = "'0.131"
$ws.Range("E40").Value = "  -7.20%  "
$ws.Range("D41").Value = "'0.321"
$ws.Range("E41").Value = "  -5.60%  "
$ws.Range("D42").Value = "'435.90"
$ws.Range("E42").Value = "  -9.94%  "
$ws.Range("D43").Value = "'48.58"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("E44").Value = "  -6.90%  "
$ws.Range("D45").Value = "'2.77"
$ws.Range("E45").Value = "  -7.37%  "
$ws.Range("D46").Value = "'8.29"
$ws.Range("E46").Value = "  -3.17%  "
$ws.Range("D48").Value = "'141.24"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").Value = "'39.57"
$ws.Range("E49").Value = "  -10.25%  "
$ws.Range("D50").Value = "2.754.77"
$ws.Range("E50").Value = "  -6.32%  "
$ws.Range("D51").Value = "'0.0345"
$ws.Range("E51").Value = "  -4.89%  "

# Reset style on the text-forced cells so no quotePrefix style artifact remains.
$textCells = @("D5","D6","D9","D10","D12","D13","D14","D19","D20","D21","D22","D23","D24","D25","D27","D29","D30","D31","D32","D33","D36","D37","D38","D39","D40","D41","D42","D43","D45","D46","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
